$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-12-06 Wednesday" "2023-12-07 Thursday"

Replace-Text "78÷9=" "51÷8="
Replace-Text "50÷4=" "11÷2="
Replace-Text "48÷3=" "68÷9="
Replace-Text "48÷8=" "74÷4="
Replace-Text "41÷6=" "51÷7="
Replace-Text "60÷7=" "96÷3="
Replace-Text "56÷3=" "27÷4="
Replace-Text "93÷9=" "12÷7="
Replace-Text "25÷2=" "81÷2="
Replace-Text "77÷5=" "29÷6="
Replace-Text "97÷6=" "87÷6="
Replace-Text "77÷2=" "26÷4="
Replace-Text "41÷7=" "67÷2="
Replace-Text "14÷2=" "12÷4="
Replace-Text "74÷6=" "83÷9="
Replace-Text "10÷8=" "42÷4="
Replace-Text "81÷5=" "47÷9="
Replace-Text "76÷4=" "21÷2="
Replace-Text "30÷3=" "12÷3="
Replace-Text "91÷3=" "60÷6="
Replace-Text "20÷9=" "35÷9="
Replace-Text "96÷4=" "10÷7="
Replace-Text "46÷7=" "28÷8="
Replace-Text "34÷8=" "75÷7="
Replace-Text "89÷6=" "30÷7="
